$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.560.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.024.91"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.23%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.47"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.06"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.95%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.545"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.50%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.76"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.67%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0858"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.496.91"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.52"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.76"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.028.22"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.64%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -14.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.599.94"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.50"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.86%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.03"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.90%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.20"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.50"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.173"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.46%  "

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.19"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.98%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.32"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.32"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.06"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0454"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.46%  "

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.50"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.46"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.40%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.286"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +10.86%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.88"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.24%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.58"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.27%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.54"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.37%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.72%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.45%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.28%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.032.79"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.325.39"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0321"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.02%  "

